$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177, shifting existing rows 177:285 down to 178:286
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new data point
$ws.Cells.Item(177, 1).Value = 4
$ws.Cells.Item(177, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(177, 3).Value = "Los Lagos"
$ws.Cells.Item(177, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(177, 5).Value = 10
$ws.Cells.Item(177, 6).Value = 100114013
$ws.Cells.Item(177, 7).Value = "Zanahoria"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 700
$ws.Cells.Item(177, 11).Value = 12500
$ws.Cells.Item(177, 12).Value = 13000
$ws.Cells.Item(177, 13).Value = 12750
$ws.Cells.Item(177, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(177, 15).Value = "Región de Ñuble"
$ws.Cells.Item(177, 16).Value = 638
$ws.Cells.Item(177, 17).Value = 20
$ws.Cells.Item(177, 18).Value = "Hortaliza"
